# Generate Report for Archive
#
# The localization status for the zh-cn / de-de handoff moved from
# "Ready for handoff" to "In Translation". That status string is shown
# in three places:
#   - Overview sheet, columns E (zh-cn) and F (de-de), row 2
#   - zh-cn sheet, column C ("Status"), row 2
#   - de-de sheet, column C ("Status"), row 2
#
# Shrinking the status text also narrows the (auto-fit) status columns
# on those same sheets, so the column widths are nudged in afterward.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Update the status values.
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Re-fit the status columns to the shorter text.
$newWidth = 13.4101845877511
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth
$zhcn.Columns.Item(3).ColumnWidth = $newWidth
$dede.Columns.Item(3).ColumnWidth = $newWidth
